$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '28.402.28'
$ws.Cells.Item(2, 5).Value = '  -0.11%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '1.819.14'
$ws.Cells.Item(3, 5).Value = '  -0.65%  '

# Row 4
$ws.Cells.Item(4, 4).NumberFormat = '@'
$ws.Cells.Item(4, 4).Value = '1.003'
$ws.Cells.Item(4, 5).Value = '  +0.16%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '314.47'
$ws.Cells.Item(5, 5).Value = '  -1.02%  '

# Row 6
$ws.Cells.Item(6, 5).Value = '  +0.14%  '

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '0.5077'
$ws.Cells.Item(7, 5).Value = '  -4.60%  '

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '0.3912'
$ws.Cells.Item(8, 5).Value = '  -3.66%  '

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '0.07684'
$ws.Cells.Item(9, 5).Value = '  +1.18%  '

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '41.82'
$ws.Cells.Item(10, 5).Value = '  -0.08%  '

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '1.105'
$ws.Cells.Item(11, 5).Value = '  -0.28%  '

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '20.89'
$ws.Cells.Item(12, 5).Value = '  -0.27%  '

# Row 13
$ws.Cells.Item(13, 2).Value = 'BinanceUSD'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '1.003'
$ws.Cells.Item(13, 5).Value = '  +0.14%  '

# Row 14
$ws.Cells.Item(14, 2).Value = 'Polkadot'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '6.232'
$ws.Cells.Item(14, 5).Value = '  -2.15%  '

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '7.489'
$ws.Cells.Item(15, 5).Value = '  -1.11%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '1.828.47'
$ws.Cells.Item(16, 5).Value = '  +0.12%  '

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '0.00001138'
$ws.Cells.Item(17, 5).Value = '  +6.00%  '

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '92.48'
$ws.Cells.Item(18, 5).Value = '  +3.50%  '

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '0.06631'
$ws.Cells.Item(19, 5).Value = '  +0.42%  '

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '17.68'
$ws.Cells.Item(20, 5).Value = '  +0.50%  '

# Row 21
$ws.Cells.Item(21, 5).Value = '  +0.04%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '6.083'
$ws.Cells.Item(22, 5).Value = '  +0.17%  '

# Row 23
$ws.Cells.Item(23, 4).Value = '28.447.36'
$ws.Cells.Item(23, 5).Value = '  -0.06%  '

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '11.22'
$ws.Cells.Item(24, 5).Value = '  -0.52%  '

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '2.258'
$ws.Cells.Item(25, 5).Value = '  +4.62%  '

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '21.08'
$ws.Cells.Item(26, 5).Value = '  +2.48%  '

# Row 27
$ws.Cells.Item(27, 2).Value = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Cells.Item(27, 4).Value = '2.035.71'
$ws.Cells.Item(27, 5).Value = '  -0.14%  '

# Row 28
$ws.Cells.Item(28, 2).Value = 'Monero'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '155.48'
$ws.Cells.Item(28, 5).Value = '  -0.90%  '

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '2.379'
$ws.Cells.Item(29, 5).Value = '  -3.91%  '

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '124.47'
$ws.Cells.Item(30, 5).Value = '  +0.56%  '

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '0.1096'
$ws.Cells.Item(31, 5).Value = '  +0.37%  '

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '1.102'
$ws.Cells.Item(32, 5).Value = '  -2.18%  '

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '5.639'
$ws.Cells.Item(33, 5).Value = '  -0.78%  '

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '3.655'
$ws.Cells.Item(34, 5).Value = '  -0.11%  '

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '0.07025'
$ws.Cells.Item(35, 5).Value = '  -2.33%  '

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '0.2204'
$ws.Cells.Item(36, 5).Value = '  -2.50%  '

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '0.02321'
$ws.Cells.Item(37, 5).Value = '  -0.88%  '

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '5.164'
$ws.Cells.Item(38, 5).Value = '  -1.09%  '

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '8.771'
$ws.Cells.Item(39, 5).Value = '  -0.17%  '

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '0.6252'
$ws.Cells.Item(40, 5).Value = '  -0.44%  '

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '11.14'
$ws.Cells.Item(41, 5).Value = '  -1.53%  '

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '1.169'
$ws.Cells.Item(42, 5).Value = '  -1.49%  '

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '1.001'
$ws.Cells.Item(43, 5).Value = '  +0.02%  '

# Row 44
$ws.Cells.Item(44, 5).Value = '  -0.55%  '

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '13.41'
$ws.Cells.Item(45, 5).Value = '  -0.33%  '

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '3.727'
$ws.Cells.Item(46, 5).Value = '  +0.57%  '

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '0.5864'
$ws.Cells.Item(47, 5).Value = '  +0.28%  '

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '123.99'
$ws.Cells.Item(48, 5).Value = '  -1.73%  '

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '1.976'
$ws.Cells.Item(49, 5).Value = '  -0.64%  '

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '1.191'
$ws.Cells.Item(50, 5).Value = '  -0.78%  '

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '0.06894'
$ws.Cells.Item(51, 5).Value = '  -0.09%  '
